$p = $ppt.ActivePresentation

# --- Slide 2 edits -------------------------------------------------------
$slide = $p.Slides.Item(2)

# "Work Structures" placeholder -> vertical anchor top (was centre) and
# retitled (singular, extra leading spaces) to "Work Structure".
$workShape = $slide.Shapes.Item(2)
$workShape.TextFrame.VerticalAnchor = 1
$workShape.TextFrame.TextRange.Text = "                   Work Structure"

# "Community Ideologies" placeholder -> retitled (singular, extra leading
# spaces) to "Community Ideology".
$communityShape = $slide.Shapes.Item(4)
$communityShape.TextFrame.TextRange.Text = "                                     Community Ideology"

# --- Presentation-level slide guide list ---------------------------------
# The canonical edit also stamps an (empty) PowerPoint 2013+ "static guides"
# extension onto the presentation root; touch the Guides collection so it
# gets persisted if the host supports it (no-op elsewhere).
try {
    $null = $p.Guides
} catch {
}
